# Apply updated crypto price/volume data to Sheet1
# (values that look numeric are apostrophe-prefixed so Excel stores them as
# text, matching the source data which keeps these as text strings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "29.219.26"
    "E2" = "  +0.31%  "
    "D3" = "1.833.23"
    "E3" = "  -0.40%  "
    "D4" = "'0.9989"
    "E4" = "  -0.09%  "
    "D5" = "'242.43"
    "E5" = "  -0.77%  "
    "D6" = "'0.6246"
    "E6" = "  -0.15%  "
    "D7" = "'1.0000"
    "E7" = "  -0.16%  "
    "D8" = "'0.07388"
    "E8" = "  -1.34%  "
    "D9" = "'0.2924"
    "E9" = "  -0.51%  "
    "D10" = "'23.22"
    "E10" = "  -0.35%  "
    "D11" = "'0.07675"
    "E11" = "  -0.51%  "
    "D12" = "1.836.19"
    "E12" = "  -1.21%  "
    "D13" = "'4.972"
    "E13" = "  -0.90%  "
    "D14" = "'0.6692"
    "E14" = "  -0.92%  "
    "D15" = "'82.62"
    "E15" = "  -0.52%  "
    "D16" = "'0.000008985"
    "E16" = "  -3.27%  "
    "D17" = "'5.889"
    "E17" = "  -1.42%  "
    "D18" = "29.188.77"
    "E18" = "  +0.17%  "
    "D19" = "2.070.79"
    "E19" = "  -2.95%  "
    "D20" = "'236.38"
    "E20" = "  +2.38%  "
    "D21" = "'12.49"
    "E21" = "  -1.71%  "
    "D22" = "'0.9995"
    "E22" = "  -0.25%  "
    "D23" = "'7.378"
    "E23" = "  +2.50%  "
    "D24" = "'0.9999"
    "E24" = "  -0.14%  "
    "D25" = "'158.29"
    "E25" = "  -1.35%  "
    "D26" = "'0.1411"
    "E26" = "  +1.20%  "
    "D27" = "'8.541"
    "E27" = "  -0.18%  "
    "D28" = "'17.69"
    "E28" = "  -1.25%  "
    "D29" = "'1.487"
    "E29" = "  -1.17%  "
    "D30" = "'0.05822"
    "E30" = "  +4.23%  "
    "D31" = "'4.107"
    "E31" = "  -1.03%  "
    "D32" = "'4.087"
    "E32" = "  -2.37%  "
    "D33" = "'1.207"
    "E33" = "  +0.03%  "
    "D34" = "'1.869"
    "E34" = "  +0.85%  "
    "D35" = "'0.7332"
    "E35" = "  -2.31%  "
    "D36" = "'1.144"
    "E36" = "  -0.01%  "
    "D37" = "'2.603"
    "E37" = "  -2.13%  "
    "D38" = "'2.852"
    "E38" = "  +3.03%  "
    "D39" = "1.224.08"
    "E39" = "  +0.16%  "
    "E40" = "  -1.52%  "
    "D41" = "'6.288"
    "E41" = "  -4.38%  "
    "D42" = "'0.9087"
    "E42" = "  +0.97%  "
    "D43" = "'1.000"
    "E43" = "  -0.08%  "
    "D44" = "'101.72"
    "E44" = "  -0.53%  "
    "D45" = "1.984.89"
    "E45" = "  -1.91%  "
    "E46" = "  -1.45%  "
    "D47" = "'0.5041"
    "E47" = "  -1.02%  "
    "B48" = "BabyDogeCoin"
    "C48" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D48" = "'0.00000000118"
    "E48" = "  -4.59%  "
    "B49" = "EnergySwap"
    "C49" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D49" = "'9.166"
    "E49" = "  +0.43%  "
    "D50" = "'0.4034"
    "E50" = "  -1.52%  "
    "D51" = "'0.1137"
    "E51" = "  +3.14%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
